$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.892.36'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.12%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.894.63'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7739'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.22%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '244.34'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.39%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3128'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.65%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.71'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.56%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07222'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.04%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08863'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +9.47%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.984.44'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.94%  '

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7717'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.427'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.27%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.46'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.14%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.190'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.47%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.957.54'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.33%  '

$ws.Range('E18').Value = '  +0.37%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.59'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.69%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007862'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.164.02'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.154'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('E24').Value = '  -0.03%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1605'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.76%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.517'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.44%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.47'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.71%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.81'
$ws.Range('D28').ClearFormats()

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.043'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.36%  '

$ws.Range('E30').Value = '  +2.11%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.561'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.09%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.543'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.113'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.46%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05494'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.92%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.249'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.42%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7512'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.71%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9993'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.714'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.60%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01956'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.64%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.785'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.17%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4506'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.87%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.97'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.25%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.092.13'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.028'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.50%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8551'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.58%  '

$ws.Range('E46').Value = '  -0.09%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.888'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.70%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.69'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.56%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.616'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.30%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.862'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.26%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.057.33'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.01%  '
